$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 207
$ws1.Range("F8").Value = 844
$ws1.Range("F9").Value = 252
$ws1.Range("F12").Value = 844
$ws1.Range("F13").Value = 291
$ws1.Range("F18").Value = 1262
$ws1.Range("F19").Value = 1200
$ws1.Range("F20").Value = 2893
$ws1.Range("F21").Value = 1447
$ws1.Range("F22").Value = 711
$ws1.Range("F24").Value = 1279
$ws1.Range("F26").Value = 1025
$ws1.Range("F28").Value = 3151
$ws1.Range("F29").Value = 612
$ws1.Range("F30").Value = 537
$ws1.Range("F31").Value = 1416

# Sheet 2: "演出" - update "想去人数" (F column) values
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F10").Value = 26

# Sheet 4: "全部类型" - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F10").Value = 207
$ws4.Range("F14").Value = 844
$ws4.Range("F15").Value = 252
$ws4.Range("F20").Value = 26
$ws4.Range("F23").Value = 844
$ws4.Range("F24").Value = 291
$ws4.Range("F29").Value = 1262
$ws4.Range("F30").Value = 1200
$ws4.Range("F31").Value = 2893
$ws4.Range("F32").Value = 1447
$ws4.Range("F33").Value = 711
$ws4.Range("F35").Value = 1279
$ws4.Range("F39").Value = 1025
$ws4.Range("F41").Value = 3151
$ws4.Range("F42").Value = 612
$ws4.Range("F43").Value = 537
$ws4.Range("F44").Value = 1416
